$wb = $excel.ActiveWorkbook

# --- Rubric-I sheet: mark several rubric items as earned ("X") ---
$wsI = $wb.Worksheets.Item("Rubric-I")
$wsI.Range("E13").Value = "X"
$wsI.Range("E15").Value = "X"
$wsI.Range("E16").Value = "X"
$wsI.Range("E17").Value = "X"

# --- Rubric-R sheet: mark every criterion row as reviewed ("x") ---
$wsR = $wb.Worksheets.Item("Rubric-R")
for ($r = 3; $r -le 16; $r++) {
    $wsR.Range("D$r").Value = "x"
}

# --- Update view/selection state to match the author's final session ---
# Rubric-I: zoomed in, selection left on B12, no longer the active tab.
$wsI.Activate() | Out-Null
$wsI.Range("B12").Select() | Out-Null
$excel.ActiveWindow.Zoom = 205

# Rubric-R: becomes the active tab, zoomed in, selection left on C10.
$wsR.Activate() | Out-Null
$wsR.Range("C10").Select() | Out-Null
$excel.ActiveWindow.Zoom = 130
